$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new measurement row at row 2 (shifting existing rows 2..21 down to 3..22),
# and drop the former last row (22) so the sheet stays 21 data rows (A1:C21).
for ($r = 21; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

$ws.Cells.Item(2, 1).Value2 = -0.6957695484161374
$ws.Cells.Item(2, 2).Value2 = 1.588029444217682
$ws.Cells.Item(2, 3).Value2 = 1.020436197519302

$ws.Range("A22:C22").ClearContents()
